$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The nested-section/nested-ligne block (rows 4-5, and the matching
# A5/A6 cell comments) is being removed from the "ligne formulaire"
# template; the former row 6 (section.lignesFormulaire) becomes row 4.

# Drop the comments that only existed for the nested-section loop (A5)
# and for the nested "ligne" loop that is being dropped (A6).
$ws.Range("A5").Comment.Delete()
$ws.Range("A6").Comment.Delete()

# A4's comment referenced the now-removed "section.sections" /
# "nestedSection" loop; replace its text with the surviving
# "section.lignesFormulaire" / "ligne" loop text (this is what used to
# live in the A6 comment).
$ws.Range("A4").Comment.Delete()
$ws.Range("A4").AddComment("Auteur:`njx:each(items=""section.lignesFormulaire"", var=""ligne"", lastCell=""J6"")")

# Remove the two now-obsolete rows (old row 4: nestedSection.ordre/name,
# old row 5: nestedLigne.code/libelle/contenu). The old row 6
# (ligne.code/libelle/contenu) shifts up to become the new row 4.
$ws.Rows("4:5").Delete()

# Update the stored selection to match the saved workbook state.
$ws.Range("C7").Select() | Out-Null
